# Juno: check in to OLPRODLOC.
# Update header row of the Chai market-trends table:
#  - Make the header labels bold (Cell.Range.Font.Bold)
#  - Capitalize "chai" -> "Chai" in several labels and tweak wording

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Header row cells that need their labels bolded (columns 2-6; column 1 "Data" is unchanged)
for ($c = 2; $c -le 6; $c++) {
    $t.Cell(1, $c).Range.Font.Bold = 1
}

# Text replacements
$d.Content.Find.Execute("Total de vendas de chai (unidades)", $true, $false, $false, $false, $false, $true, 1, $false, "Total de vendas de Chai (unidades)", 2)
$d.Content.Find.Execute("Vendas de chai artesanal (unidades)", $true, $false, $false, $false, $false, $true, 1, $false, "Vendas de Chai artesanal (unidades)", 2)
$d.Content.Find.Execute("Vendas de chai pronto (unidades)", $true, $false, $false, $false, $false, $true, 1, $false, "Vendas de Chai pré-fabricado (unidades)", 2)
$d.Content.Find.Execute("Participação nas redes sociais (exibições)", $true, $false, $false, $false, $false, $true, 1, $false, "Engajamento nas redes sociais (visualizações)", 2)
$d.Content.Find.Execute("Pesquisas online por chai", $true, $false, $false, $false, $false, $true, 1, $false, "Pesquisas online por Chai", 2)
